$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 header text change: "身分（輸入數字）" -> "身分（教授、副教授、助理教授或博士後研究員）"
$ws.Range("D1").Value = "身分（教授、副教授、助理教授或博士後研究員）"

# Column D is widened to fit the longer header text
$ws.Columns.Item(4).ColumnWidth = 46.15

# View changes: zoom level reduced, and the active selection moved
$excel.ActiveWindow.Zoom = 145
$ws.Range("D6").Select()
